$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.423.26"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "2.658.08"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.90%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("D9").Value = "2.656.22"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.11%  "

$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.42%  "

$ws.Range("D14").Value = "3.144.74"
$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000185"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.20%  "

$ws.Range("D16").Value = "72.350.72"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.88%  "

$ws.Range("D18").Value = "2.666.23"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("E19").Value = "  +3.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.87%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.00%  "

$ws.Range("E25").Value = "  -0.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.04%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").Value = "0.0₃0955"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "498.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.84%  "

$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.91"
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = "  -2.84%  "

$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("E42").Value = "  -6.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.67%  "

$ws.Range("E44").Value = "  -3.48%  "

$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.70%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.551"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.18%  "

$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.602"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
